$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated sensor-reading table (A:timestamp, B:label, C:ax, D:ay, E:az,
# F:gx, G:gy, H:gz). Rows 2-11 receive new readings for the existing
# timestamps 0-900; rows 12-31 hold what used to be rows 2-21 (shifted
# down by 10), with 10 brand new rows (timestamps 2000-2900) appended.
$numRows = 30
$numCols = 8
$data = New-Object 'object[,]' $numRows, $numCols

$data[0,0] = 0
$data[0,1] = "falling"
$data[0,2] = -0.3973007202148437
$data[0,3] = 1.798293828964233
$data[0,4] = 1.168385148048401
$data[0,5] = -0.0282525178045034
$data[0,6] = 0.0256563406437635
$data[0,7] = 0.06856962293386459
$data[1,0] = 100
$data[1,1] = "falling"
$data[1,2] = -0.8141142129898071
$data[1,3] = 1.826734185218811
$data[1,4] = 1.343665383756161
$data[1,5] = 0.0575740486383438
$data[1,6] = -0.064446285367012
$data[1,7] = 0.1545489132404327
$data[2,0] = 200
$data[2,1] = "falling"
$data[2,2] = -0.9248467683792115
$data[2,3] = 1.964664489030838
$data[2,4] = 1.209049716591835
$data[2,5] = -0.0140499006956815
$data[2,6] = 0.0291688162833452
$data[2,7] = 0.0418442711234092
$data[3,0] = 300
$data[3,1] = "falling"
$data[3,2] = -0.7549184560775757
$data[3,3] = 1.929333925247193
$data[3,4] = 1.357575602829456
$data[3,5] = -0.0314595587551593
$data[3,6] = -0.0210748501121997
$data[3,7] = 0.0074830991216003
$data[4,0] = 400
$data[4,1] = "falling"
$data[4,2] = -0.4957029819488525
$data[4,3] = 1.878820419311524
$data[4,4] = 1.164325326681137
$data[4,5] = 0.0100792767480015
$data[4,6] = -0.0314595587551593
$data[4,7] = -0.0143553335219621
$data[5,0] = 500
$data[5,1] = "falling"
$data[5,2] = -0.8541634678840635
$data[5,3] = 1.869282335042953
$data[5,4] = 1.307794235646725
$data[5,5] = -0.00137444678694
$data[5,6] = 0.0216857157647609
$data[5,7] = -0.016951510682702
$data[6,0] = 600
$data[6,1] = "falling"
$data[6,2] = -0.8955824375152587
$data[6,3] = 1.77057421207428
$data[6,4] = 1.258034527301789
$data[6,5] = 0.0010690141934901
$data[6,6] = 0.0012217304902151
$data[6,7] = -0.0684169083833694
$data[7,0] = 700
$data[7,1] = "falling"
$data[7,2] = -0.5728458166122438
$data[7,3] = 1.762963086366654
$data[7,4] = 1.322432711720467
$data[7,5] = -0.0058032199740409
$data[7,6] = -0.0200058370828628
$data[7,7] = -0.07605272531509399
$data[8,0] = 800
$data[8,1] = "falling"
$data[8,2] = -0.5884580612182617
$data[8,3] = 1.777032017707825
$data[8,4] = 1.320214748382568
$data[8,5] = -0.0381790772080421
$data[8,6] = -0.0108428578823804
$data[8,7] = 0.0048869219608604
$data[9,0] = 900
$data[9,1] = "falling"
$data[9,2] = -0.7303044199943542
$data[9,3] = 1.775961980223656
$data[9,4] = 1.753339484333992
$data[9,5] = -0.0554360225796699
$data[9,6] = 0.0329867228865623
$data[9,7] = 0.012980886735022
$data[10,0] = 1000
$data[10,1] = "falling"
$data[10,2] = -0.7120251655578613
$data[10,3] = 1.784507364034653
$data[10,4] = 2.399995267391204
$data[10,5] = 0.0018325957935303
$data[10,6] = 0.1020144969224929
$data[10,7] = 0.0574213340878486
$data[11,0] = 1100
$data[11,1] = "falling"
$data[11,2] = -0.8870223164558411
$data[11,3] = 1.693844005465508
$data[11,4] = 2.417137637734413
$data[11,5] = -0.1429424732923507
$data[11,6] = 0.204487144947052
$data[11,7] = 0.0444404482841491
$data[12,0] = 1200
$data[12,1] = "falling"
$data[12,2] = -1.0569589138031
$data[12,3] = 1.552494168281555
$data[12,4] = 2.973462641239166
$data[12,5] = -0.2125810980796814
$data[12,6] = 0.4230241775512695
$data[12,7] = 0.0167987942695617
$data[13,0] = 1300
$data[13,1] = "falling"
$data[13,2] = -1.252092391252517
$data[13,3] = 1.304344907402993
$data[13,4] = 3.757617935538292
$data[13,5] = -0.0858265683054924
$data[13,6] = 0.5377141237258911
$data[13,7] = -0.1299615800380706
$data[14,0] = 1400
$data[14,1] = "falling"
$data[14,2] = -2.462630152702332
$data[14,3] = 0.6072362959384919
$data[14,4] = 3.575700670480729
$data[14,5] = -0.1346957832574844
$data[14,6] = 0.5236642360687256
$data[14,7] = 0.1327104717493057
$data[15,0] = 1500
$data[15,1] = "falling"
$data[15,2] = -2.071417987346649
$data[15,3] = 0.361852638423443
$data[15,4] = 3.242989569902419
$data[15,5] = -0.0175623763352632
$data[15,6] = -0.0714712366461753
$data[15,7] = -0.2777909636497497
$data[16,0] = 1600
$data[16,1] = "falling"
$data[16,2] = -3.716687679290772
$data[16,3] = 0.0018857717514038
$data[16,4] = 4.452265739440918
$data[16,5] = 0.0302378293126821
$data[16,6] = 0.5580254197120667
$data[16,7] = -0.6604980230331421
$data[17,0] = 1700
$data[17,1] = "falling"
$data[17,2] = -5.890828639268875
$data[17,3] = 0.2718941420316697
$data[17,4] = 4.244596153497697
$data[17,5] = 0.4694499373435974
$data[17,6] = -1.051451802253723
$data[17,7] = -0.956462264060974
$data[18,0] = 1800
$data[18,1] = "falling"
$data[18,2] = -7.35106348991394
$data[18,3] = 2.290824487805367
$data[18,4] = 0.9694806933402997
$data[18,5] = -1.834733724594116
$data[18,6] = 0.0054977871477603
$data[18,7] = -1.47539234161377
$data[19,0] = 1900
$data[19,1] = "falling"
$data[19,2] = -12.02078425884247
$data[19,3] = 0.3663914650678617
$data[19,4] = 3.381503492593767
$data[19,5] = 2.41551399230957
$data[19,6] = 0.4040873646736145
$data[19,7] = 0.6867652535438538
$data[20,0] = 2000
$data[20,1] = "falling"
$data[20,2] = -19.42419934272767
$data[20,3] = 8.062076985836041
$data[20,4] = 6.704323172569278
$data[20,5] = -0.168751522898674
$data[20,6] = 0.6850853562355042
$data[20,7] = 0.1965458989143371
$data[21,0] = 2100
$data[21,1] = "falling"
$data[21,2] = -1.815191864967314
$data[21,3] = 4.267351135611525
$data[21,4] = 2.374830707907669
$data[21,5] = -0.2154827117919922
$data[21,6] = 0.3124575614929199
$data[21,7] = 0.2249511331319809
$data[22,0] = 2200
$data[22,1] = "falling"
$data[22,2] = 3.242185950279228
$data[22,3] = 1.317612782120705
$data[22,4] = 0.6843594610691084
$data[22,5] = -0.1125519201159477
$data[22,6] = -0.1751656085252761
$data[22,7] = 0.4828889667987823
$data[23,0] = 2300
$data[23,1] = "falling"
$data[23,2] = -0.3276352286338837
$data[23,3] = 2.805960930883888
$data[23,4] = 2.484445497393612
$data[23,5] = -0.1719585657119751
$data[23,6] = -0.7235698699951172
$data[23,7] = -0.6099489331245422
$data[24,0] = 2400
$data[24,1] = "falling"
$data[24,2] = -0.4897050857543918
$data[24,3] = 1.077956080436697
$data[24,4] = 1.630782932043069
$data[24,5] = -0.0491746515035629
$data[24,6] = -0.7188356518745422
$data[24,7] = 0.3561344444751739
$data[25,0] = 2500
$data[25,1] = "falling"
$data[25,2] = 0.7483568191528327
$data[25,3] = -0.1515689864754661
$data[25,4] = 0.4959011934697636
$data[25,5] = -0.0861319974064827
$data[25,6] = 0.2157881408929824
$data[25,7] = 0.3926336467266083
$data[26,0] = 2600
$data[26,1] = "falling"
$data[26,2] = 1.010472297668458
$data[26,3] = 2.143007203936576
$data[26,4] = 1.17438416928053
$data[26,5] = -0.0914770737290382
$data[26,6] = 0.1360702365636825
$data[26,7] = -0.0354301854968071
$data[27,0] = 2700
$data[27,1] = "falling"
$data[27,2] = 0.2876673340797431
$data[27,3] = 2.055301316082478
$data[27,4] = 0.5345815420150758
$data[27,5] = 0.1240056455135345
$data[27,6] = 0.0096211275085806
$data[27,7] = -0.4243986308574676
$data[28,0] = 2800
$data[28,1] = "falling"
$data[28,2] = -0.02689480781555212
$data[28,3] = 1.347976356744766
$data[28,4] = 0.6012542694807046
$data[28,5] = 0.009468411095440299
$data[28,6] = 0.0371100641787052
$data[28,7] = 0.0630718395113945
$data[29,0] = 2900
$data[29,1] = "falling"
$data[29,2] = 0.3917713761329649
$data[29,3] = 1.835358053445816
$data[29,4] = 1.164814613759517
$data[29,5] = 0.0762054398655891
$data[29,6] = 0.1154535338282585
$data[29,7] = 0.0804814994335174

$range = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(1 + $numRows, $numCols))
$range.Value2 = $data

$wb.Save()
